## Figure 1.1 caption tweaks:
##   a. Original memory arrangement
##   b. Refactored and reduced data
##   c. Initial storage layout
##   d. Storage layout after system optimization
## All four captions become italic, and the two captions that live
## inside the little two-box diagram group get wider text boxes so the
## new text fits.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# ---------------------------------------------------------------------
# (c) "Initial storage layout" -> "c. Initial storage layout"
#     TextBox 151 (id 152), inside Group 138 (id 139)
# ---------------------------------------------------------------------
$grpC = $s.Shapes.Item(15)
$tbC = $grpC.GroupItems.Item(12)

$subC = $tbC.TextFrame.TextRange.Characters(1, 7)   # "Initial"
$subC.Text = "c. Initial"

$tbC.TextFrame.TextRange.Font.Italic = $true
$tbC.Width = 120.2702   # 1527431 EMU (was 1358709 EMU)

# ---------------------------------------------------------------------
# (d) "Layout after system optimization" ->
#     "d. Storage layout after system optimization"
#     TextBox 163 (id 164), inside Group 152 (id 153)
# ---------------------------------------------------------------------
$grpD = $s.Shapes.Item(16)
$tbD = $grpD.GroupItems.Item(10)

$subD = $tbD.TextFrame.TextRange.Characters(1, 7)   # "Layout "
$subD.Text = "d. Storage layout "

$tbD.TextFrame.TextRange.Font.Italic = $true
$tbD.Width = 211.0846   # 2680774 EMU (was 2067474 EMU)

# The wider caption nudges the enclosing group's tight bounding box by
# one EMU (PowerPoint recomputes grpSpPr/xfrm/off on resize).
$grpD.Left = 81.8722   # 1039776 EMU (was 1039777 EMU)

# ---------------------------------------------------------------------
# (b) "Refactored and reduced data" -> "b. Refactored and reduced data"
#     TextBox 55 (id 56)
# ---------------------------------------------------------------------
$tbB = $s.Shapes.Item(25)

$subB = $tbB.TextFrame.TextRange.Characters(1, 15)  # "Refactored and "
$subB.Text = "b. Refactored and "

$tbB.TextFrame.TextRange.Font.Italic = $true

# ---------------------------------------------------------------------
# (a) "Original memory arrangement" -> "a. Original memory arrangement"
#     TextBox 60 (id 61)
# ---------------------------------------------------------------------
$tbA = $s.Shapes.Item(29)

$subA = $tbA.TextFrame.TextRange.Characters(1, $tbA.TextFrame.TextRange.Length)
$subA.Text = "a. Original memory arrangement"

$tbA.TextFrame.TextRange.Font.Italic = $true
